$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture the current (pre-edit) row 3 content, which will move down to row 5 ---
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2
$i3 = $ws.Range("I3").Value2
$j3 = $ws.Range("J3").Value2
$k3 = $ws.Range("K3").Value2
$l3 = $ws.Range("L3").Value2
$m3 = $ws.Range("M3").Value2
$n3 = $ws.Range("N3").Value2
$q3 = $ws.Range("Q3").Value2

# --- Step 2: capture the current (pre-edit) row 4 content, which will move up to row 3 ---
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$n4 = $ws.Range("N4").Value2

# --- Step 3: give row 5 the same cell formatting row 3 currently has, then fill in its values ---
$ws.Range("A3:Q3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A5").Value2 = $a3
$ws.Range("B5").Value2 = $b3
$ws.Range("C5").Value2 = $c3
$ws.Range("D5").Formula = "=IFERROR(LEFT(C5,FIND("" "",C5)-1),C5)"
$ws.Range("E5").Value2 = $e3
$ws.Range("F5").Value2 = $f3
$ws.Range("G5").Value2 = $g3
$ws.Range("H5").Value2 = $h3
$ws.Range("I5").Value2 = $i3
$ws.Range("J5").Value2 = $j3
$ws.Range("K5").Value2 = $k3
$ws.Range("L5").Value2 = $l3
$ws.Range("M5").Value2 = $m3
$ws.Range("N5").Value2 = $n3
$ws.Range("O5").Formula = "=IF(N5=Data!`$B`$18,""Programming"","""")"
$ws.Range("P5").Formula = "=IF(N5=Data!`$B`$19,""9"","""")"
$ws.Range("Q5").Value2 = $q3

# --- Step 4: overwrite row 3 with the old row-4 data (only A, B, N actually differ) ---
$ws.Range("A3").Value2 = $a4
$ws.Range("B3").Value2 = $b4
$ws.Range("N3").Value2 = $n4

# --- Step 5: row 4 no longer holds any data; wipe it completely ---
$ws.Range("A4:Q4").Clear()

# --- Step 6: tidy up the two plain (non-x14) data validations so they skip row 4 ---
$ws.Range("C2:C4").Validation.Delete()
$ws.Range("C2:C5").Validation.Add(3, 1, 1, "Category")
$ws.Range("C4").Validation.Delete()

$ws.Range("E2:E4").Validation.Delete()
$ws.Range("E2:E5").Validation.Add(3, 1, 1, "=INDIRECT(SUBSTITUTE(D2,"" "",""""))")
$ws.Range("E4").Validation.Delete()

# --- Step 7: selection / column widths to mirror the manual edit ---
$ws.Columns("A").AutoFit()
$ws.Columns("B").AutoFit()
$ws.Range("A5:XFD5").Select()
